$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-18 23:48:24"
$ws.Range("I2").Value = "2.1 mm"
$ws.Range("E3").Value = "2026-02-18 23:48:27"
$ws.Range("H3").Value = "'73%"
$ws.Range("I3").Value = "0.8 mm"
$ws.Range("E4").Value = "2026-02-18 23:48:30"
$ws.Range("H4").Value = "'71%"
$ws.Range("J4").Value = "1012.0 hPa"
$ws.Range("E5").Value = "2026-02-18 23:48:32"
$ws.Range("I5").Value = "2.2 mm"
$ws.Range("N5").Value = "-4.8 °C 23:27 TU"
$ws.Range("O5").Value = "0.3 °C"
$ws.Range("E6").Value = "2026-02-18 23:48:35"
$ws.Range("J6").Value = "1011.7 hPa"
$ws.Range("E7").Value = "2026-02-18 23:48:37"
$ws.Range("J7").Value = "1013.3 hPa"
$ws.Range("L7").Value = "46.1 km/h - 290º 23:11 TU"
$ws.Range("E8").Value = "2026-02-18 23:48:40"
$ws.Range("J8").Value = "1013.1 hPa"
$ws.Range("E9").Value = "2026-02-18 23:48:43"
$ws.Range("E10").Value = "2026-02-18 23:48:45"
$ws.Range("H10").Value = "'83%"
$ws.Range("E11").Value = "2026-02-18 23:48:47"
$ws.Range("E12").Value = "2026-02-18 23:48:49"
$ws.Range("E13").Value = "2026-02-18 23:48:52"
$ws.Range("J13").Value = "1014.3 hPa"
$ws.Range("E14").Value = "2026-02-18 23:48:54"
$ws.Range("N14").Value = "8.0 °C 23:19 TU"
$ws.Range("O14").Value = "12.1 °C"
$ws.Range("E15").Value = "2026-02-18 23:48:57"
$ws.Range("H15").Value = "'78%"
$ws.Range("O15").Value = "10.7 °C"
$ws.Range("E16").Value = "2026-02-18 23:48:59"
$ws.Range("H16").Value = "'55%"
$ws.Range("I16").Value = "2.8 mm"
$ws.Range("N16").Value = "-4.8 °C 23:29 TU"
$ws.Range("E17").Value = "2026-02-18 23:49:02"
$ws.Range("H17").Value = "'87%"
$ws.Range("N17").Value = "1.5 °C 23:21 TU"
$ws.Range("E18").Value = "2026-02-18 23:49:04"
$ws.Range("H18").Value = "'77%"
$ws.Range("J18").Value = "1012.2 hPa"
$ws.Range("E19").Value = "2026-02-18 23:49:06"
$ws.Range("E20").Value = "2026-02-18 23:49:09"
$ws.Range("H20").Value = "'78%"
$ws.Range("N20").Value = "-4.7 °C 23:17 TU"
$ws.Range("O20").Value = "-0.7 °C"
$ws.Range("E21").Value = "2026-02-18 23:49:11"
$ws.Range("H21").Value = "'71%"
$ws.Range("J21").Value = "1013.8 hPa"
$ws.Range("O21").Value = "6.6 °C"
$ws.Range("E22").Value = "2026-02-18 23:49:14"
$ws.Range("O22").Value = "-2.0 °C"
$ws.Range("E23").Value = "2026-02-18 23:49:16"
$ws.Range("H23").Value = "'61%"
$ws.Range("I23").Value = "0.9 mm"
$ws.Range("N23").Value = "-5.4 °C 23:23 TU"
$ws.Range("O23").Value = "-0.3 °C"
$ws.Range("E24").Value = "2026-02-18 23:49:19"
$ws.Range("J24").Value = "1014.0 hPa"
$ws.Range("E25").Value = "2026-02-18 23:49:21"
$ws.Range("O25").Value = "1.5 °C"
$ws.Range("E26").Value = "2026-02-18 23:49:24"
$ws.Range("J26").Value = "1011.3 hPa"
$ws.Range("O26").Value = "5.4 °C"
$ws.Range("E27").Value = "2026-02-18 23:49:26"
$ws.Range("H27").Value = "'61%"
$ws.Range("N27").Value = "-2.9 °C 23:24 TU"
$ws.Range("O27").Value = "1.2 °C"
$ws.Range("E28").Value = "2026-02-18 23:49:29"
$ws.Range("H28").Value = "'71%"
$ws.Range("J28").Value = "1011.9 hPa"
$ws.Range("O28").Value = "10.0 °C"
$ws.Range("E29").Value = "2026-02-18 23:49:31"
$ws.Range("E30").Value = "2026-02-18 23:49:34"
$ws.Range("J30").Value = "1011.5 hPa"
$ws.Range("E31").Value = "2026-02-18 23:49:36"
$ws.Range("H31").Value = "'72%"
$ws.Range("J31").Value = "1010.4 hPa"
$ws.Range("E32").Value = "2026-02-18 23:49:39"
$ws.Range("E33").Value = "2026-02-18 23:49:41"
$ws.Range("J33").Value = "1013.1 hPa"
$ws.Range("O33").Value = "4.9 °C"
$ws.Range("E34").Value = "2026-02-18 23:49:44"
$ws.Range("E35").Value = "2026-02-18 23:49:46"
$ws.Range("I35").Value = "2.6 mm"
$ws.Range("J35").Value = "1013.8 hPa"
$ws.Range("L35").Value = "75.2 km/h - 252º 23:15 TU"
$ws.Range("E36").Value = "2026-02-18 23:49:49"
$ws.Range("J36").Value = "1012.0 hPa"
$ws.Range("O36").Value = "12.0 °C"
$ws.Range("E37").Value = "2026-02-18 23:49:51"
$ws.Range("J37").Value = "1013.5 hPa"
$ws.Range("E38").Value = "2026-02-18 23:49:54"
$ws.Range("H38").Value = "'74%"
$ws.Range("E39").Value = "2026-02-18 23:49:56"
$ws.Range("H39").Value = "'45%"
$ws.Range("I39").Value = "0.6 mm"
$ws.Range("N39").Value = "-2.5 °C 23:29 TU"
$ws.Range("E40").Value = "2026-02-18 23:49:59"
$ws.Range("J40").Value = "1014.5 hPa"
$ws.Range("E41").Value = "2026-02-18 23:50:01"
$ws.Range("J41").Value = "1013.6 hPa"
$ws.Range("K41").Value = "11.9 MJ/m2"
$ws.Range("L41").Value = "37.1 km/h - 252º 23:06 TU"
$ws.Range("O41").Value = "11.4 °C"
$ws.Range("E42").Value = "2026-02-18 23:50:04"
$ws.Range("H42").Value = "'83%"
$ws.Range("E43").Value = "2026-02-18 23:50:06"
$ws.Range("H43").Value = "'77%"
$ws.Range("E44").Value = "2026-02-18 23:50:09"
$ws.Range("H44").Value = "'76%"
$ws.Range("I44").Value = "0.3 mm"
$ws.Range("O44").Value = "-1.7 °C"
$ws.Range("E45").Value = "2026-02-18 23:50:11"
$ws.Range("H45").Value = "'67%"
$ws.Range("I45").Value = "1.7 mm"
$ws.Range("J45").Value = "1011.3 hPa"
$ws.Range("O45").Value = "7.0 °C"
$ws.Range("E46").Value = "2026-02-18 23:50:14"
$ws.Range("H46").Value = "'83%"
$ws.Range("J46").Value = "1014.1 hPa"
$ws.Range("L46").Value = "31.3 km/h - 304º 23:05 TU"
$ws.Range("O46").Value = "11.0 °C"
